# Scheduled-runner data refresh: overwrite the static market-price /
# profit columns (H:N -> currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ,
# LeveProfitHQ) on the affected rows of each job sheet with freshly
# pulled values. These are plain literals (no formulas in this workbook),
# so each cell is just re-written with its new value.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 970.5
$ws.Range("I98").Value = 863.125
$ws.Range("J98").Value = 1400
$ws.Range("K98").Value = 863.125
$ws.Range("L98").Value = 1400
$ws.Range("M98").Value = 634.875
$ws.Range("N98").Value = -4396
$ws.Range("H122").Value = 970.5
$ws.Range("I122").Value = 863.125
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 2589.375
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -139.375
$ws.Range("N122").Value = -9100
$ws.Range("H138").Value = 1595.5526
$ws.Range("J138").Value = 2503.4
$ws.Range("L138").Value = 7510.200000000001
$ws.Range("N138").Value = -17790.2
$ws.Range("H139").Value = 47540
$ws.Range("J139").Value = 47540
$ws.Range("L139").Value = 47540
$ws.Range("N139").Value = -57820
$ws.Range("H140").Value = 62500
$ws.Range("J140").Value = 62500
$ws.Range("L140").Value = 62500
$ws.Range("N140").Value = -72860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8773308
$ws.Range("I61").Value = 11364842
$ws.Range("J61").Value = 1963.6154
$ws.Range("K61").Value = 11364842
$ws.Range("L61").Value = 1963.6154
$ws.Range("M61").Value = -11364630
$ws.Range("N61").Value = -2387.6154
$ws.Range("H74").Value = 11629520
$ws.Range("I74").Value = 14287183
$ws.Range("J74").Value = 2241
$ws.Range("K74").Value = 14287183
$ws.Range("L74").Value = 2241
$ws.Range("M74").Value = -14286309
$ws.Range("N74").Value = -3989
$ws.Range("H77").Value = 11629520
$ws.Range("I77").Value = 14287183
$ws.Range("J77").Value = 2241
$ws.Range("K77").Value = 71435915
$ws.Range("L77").Value = 11205
$ws.Range("M77").Value = -71431547
$ws.Range("N77").Value = -19941
$ws.Range("H97").Value = 7841.706
$ws.Range("I97").Value = 8177.615
$ws.Range("J97").Value = 6750
$ws.Range("K97").Value = 8177.615
$ws.Range("L97").Value = 6750
$ws.Range("M97").Value = -7681.615
$ws.Range("N97").Value = -7742
$ws.Range("H122").Value = 2962.2258
$ws.Range("I122").Value = 3186.9
$ws.Range("J122").Value = 2026.0834
$ws.Range("K122").Value = 9560.700000000001
$ws.Range("L122").Value = 6078.2502
$ws.Range("M122").Value = -7110.700000000001
$ws.Range("N122").Value = -10978.2502
$ws.Range("H132").Value = 3473737.5
$ws.Range("I132").Value = 4387264
$ws.Range("J132").Value = 2336.5334
$ws.Range("K132").Value = 13161792
$ws.Range("L132").Value = 7009.600199999999
$ws.Range("M132").Value = -13159262
$ws.Range("N132").Value = -12069.6002
$ws.Range("H136").Value = 8773308
$ws.Range("I136").Value = 11364842
$ws.Range("J136").Value = 1963.6154
$ws.Range("K136").Value = 34094526
$ws.Range("L136").Value = 5890.8462
$ws.Range("M136").Value = -34091976
$ws.Range("N136").Value = -10990.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3215.75
$ws.Range("I134").Value = 2368.282
$ws.Range("K134").Value = 7104.846
$ws.Range("M134").Value = -4569.846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 118
$ws.Range("I7").Value = 125
$ws.Range("J7").Value = 111
$ws.Range("K7").Value = 125
$ws.Range("L7").Value = 111
$ws.Range("M7").Value = -12
$ws.Range("N7").Value = -337
$ws.Range("H31").Value = 5467930.5
$ws.Range("I31").Value = 3597.7917
$ws.Range("J31").Value = 25643928
$ws.Range("K31").Value = 3597.7917
$ws.Range("L31").Value = 25643928
$ws.Range("M31").Value = -3302.7917
$ws.Range("N31").Value = -25644518
$ws.Range("H34").Value = 5467930.5
$ws.Range("I34").Value = 3597.7917
$ws.Range("J34").Value = 25643928
$ws.Range("K34").Value = 3597.7917
$ws.Range("L34").Value = 25643928
$ws.Range("M34").Value = -3395.7917
$ws.Range("N34").Value = -25644332
$ws.Range("H107").Value = 760.4286
$ws.Range("I107").Value = 771.26666
$ws.Range("J107").Value = 733.3333
$ws.Range("K107").Value = 771.26666
$ws.Range("L107").Value = 733.3333
$ws.Range("M107").Value = 1148.73334
$ws.Range("N107").Value = -4573.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1595.6111
$ws.Range("I113").Value = 495
$ws.Range("J113").Value = 1910.0714
$ws.Range("K113").Value = 1485
$ws.Range("L113").Value = 5730.2142
$ws.Range("M113").Value = 685
$ws.Range("N113").Value = -10070.2142
$ws.Range("H118").Value = 1966.6666
$ws.Range("I118").Value = 1000
$ws.Range("J118").Value = 2450
$ws.Range("K118").Value = 3000
$ws.Range("L118").Value = 7350
$ws.Range("M118").Value = -1757
$ws.Range("N118").Value = -9836
$ws.Range("H122").Value = 817.65515
$ws.Range("I122").Value = 806.7368
$ws.Range("J122").Value = 838.4
$ws.Range("K122").Value = 7260.6312
$ws.Range("L122").Value = 7545.599999999999
$ws.Range("M122").Value = -4810.6312
$ws.Range("N122").Value = -12445.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4031.65
$ws.Range("I132").Value = 2944.2856
$ws.Range("J132").Value = 6568.8335
$ws.Range("K132").Value = 8832.856800000001
$ws.Range("L132").Value = 19706.5005
$ws.Range("M132").Value = -6302.856800000001
$ws.Range("N132").Value = -24766.5005
$ws.Range("H138").Value = 52800
$ws.Range("J138").Value = 52800
$ws.Range("L138").Value = 52800
$ws.Range("N138").Value = -63080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3829.8628
$ws.Range("I7").Value = 3545.7812
$ws.Range("J7").Value = 4308.316
$ws.Range("K7").Value = 3545.7812
$ws.Range("L7").Value = 4308.316
$ws.Range("M7").Value = -3433.7812
$ws.Range("N7").Value = -4532.316
$ws.Range("H16").Value = 1740.0714
$ws.Range("I16").Value = 1740.0714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1740.0714
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1570.0714
# N16 no longer carries a profit figure for this row - remove it outright
# (matches the source row losing its HQ-profit cell entirely).
$ws.Range("N16").ClearContents()
$ws.Range("H68").Value = 1691.2174
$ws.Range("I68").Value = 1722
$ws.Range("J68").Value = 1643.3334
$ws.Range("K68").Value = 1722
$ws.Range("L68").Value = 1643.3334
$ws.Range("M68").Value = -973
$ws.Range("N68").Value = -3141.3334
$ws.Range("H71").Value = 1691.2174
$ws.Range("I71").Value = 1722
$ws.Range("J71").Value = 1643.3334
$ws.Range("K71").Value = 8610
$ws.Range("L71").Value = 8216.666999999999
$ws.Range("M71").Value = -4866
$ws.Range("N71").Value = -15704.667
$ws.Range("H126").Value = 3829.8628
$ws.Range("I126").Value = 3545.7812
$ws.Range("J126").Value = 4308.316
$ws.Range("K126").Value = 10637.3436
$ws.Range("L126").Value = 12924.948
$ws.Range("M126").Value = -8167.3436
$ws.Range("N126").Value = -17864.948
$ws.Range("H136").Value = 16671599
$ws.Range("I136").Value = 25001488
$ws.Range("J136").Value = 11820.5
$ws.Range("K136").Value = 75004464
$ws.Range("L136").Value = 35461.5
$ws.Range("M136").Value = -75001914
$ws.Range("N136").Value = -40561.5
$ws.Range("H139").Value = 51700
$ws.Range("J139").Value = 51700
$ws.Range("L139").Value = 51700
$ws.Range("N139").Value = -61980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 795.3
$ws.Range("I107").Value = 960.73334
$ws.Range("K107").Value = 2882.20002
$ws.Range("M107").Value = -962.2000200000002
$ws.Range("H126").Value = 1795.0605
$ws.Range("I126").Value = 1125.2413
$ws.Range("J126").Value = 6651.25
$ws.Range("K126").Value = 3375.7239
$ws.Range("L126").Value = 19953.75
$ws.Range("M126").Value = -905.7239
$ws.Range("N126").Value = -24893.75
$ws.Range("H136").Value = 1361.6
$ws.Range("I136").Value = 1043.3334
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 3130.0002
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -580.0001999999999
$ws.Range("N136").Value = -32100
$ws.Range("H138").Value = 44771.6
$ws.Range("J138").Value = 44771.6
$ws.Range("L138").Value = 44771.6
$ws.Range("N138").Value = -55051.6
